$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDF_CH_sub_1_stratification")

# Row 2 (patient 31192): race American Indian or Alaska Native -> White
$ws.Range("C2").Value = "White"

# Row 3 (patient 34485): race Other Race -> White
$ws.Range("C3").Value = "White"

# Row 4 (patient 58288): ethnicity Not Hispanic or Latino -> Hispanic or Latino
#                         race White -> Native Hawaiian or Other Pacific Islander
$ws.Range("B4").Value = "Hispanic or Latino"
$ws.Range("C4").Value = "Native Hawaiian or Other Pacific Islander"

# Row 5 (patient 86683): race Native Hawaiian or Other Pacific Islander -> Other Race
$ws.Range("C5").Value = "Other Race"
